$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked_lbl" and "is_enabled_lbl" header columns (D1:E1).
# The remaining columns (order_by, rem) shift left to take their place,
# leaving the row ending at E1 (old F1/G1 become empty).
$ws.Range("D1:E1").Delete(-4159)
